# UniformF-HW30.xlsx -- "New simulation files for schemes report"
#
# The original sheet1 had:
#   - Row 1: B1:AD1 = 0..28 (29 sequential indices)
#   - Row 2: A2=0, B2="HKL", C2:T2 and U2:AD2 = two back-to-back copies
#            of the same 18 "HKL"-style labels
#   - Rows 3-19: A=index, B="HKL" name (17 names), C:T = 1
#
# The new sheet1 has:
#   - Row 1: B1:T1 = 0..18 (19 sequential indices; U1:AD1 removed)
#   - Row 2: A2=0, B2="HKL", C2:T2 = 18 "HKL"-style labels, a DIFFERENT
#            reordered/renamed set (U2:AD2 removed)
#   - Rows 3-29: A=index, B="HKL" name (27 renamed/reordered names),
#            C:T = 1  (27 data rows -- 10 new rows appended)
#
# Net effect: nearly every string value changed (renamed/reordered),
# 10 rows were appended, and the stale U:AD "echo" columns were dropped.
#
# Column A (rows 2-19) already carries the bordered/bold/centred style
# used throughout ("s=1" in the OOXML) and is left alone; only its
# numeric values get refreshed/extended. Columns B:AD for rows 1-19 are
# cleared outright (wiping the old shared-string pool along with them so
# it gets rebuilt from scratch, cleanly, in the exact order the new
# content is written) and then repopulated with the target values. The
# header style is restored onto B1:T1 -- and onto the ten freshly-added
# A20:A29 cells -- by copying the still-intact format from A2, which
# reproduces the existing style index exactly instead of growing the
# stylesheet with new/orphaned entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Wipe the string-bearing region (values + formatting) -----------
#        This also drops the obsolete U:AD "echo" columns and resets the
#        shared-string pool to empty so it rebuilds cleanly below.
$ws.Range("B1:AD19").Clear()

# --- 2. Restore the header style (s=1) onto B1:T1 and the new A20:A29 --
#        rows by copying format from A2, which still carries it.
$fmtSrc = $ws.Range("A2")
$fmtSrc.Copy()
$ws.Range("B1:T1").PasteSpecial(-4122)
$ws.Range("A20:A29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Header row (row 1): B1:T1 = 0..18 -------------------------------
for ($i = 0; $i -le 18; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $i
}

# --- 4. The full ordered list of "HKL" labels used down column B -------
#        (row 2 = "HKL" itself, rows 3-29 = the 27 scheme names)
$hklNames = @(
    "HKL",
    "Spiral5",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Thomas Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

# --- 5. Row 2: A2=0, B2="HKL" -------------------------------------------
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = $hklNames[0]

# --- 6. Rows 3..29: A=index (1..27), B=scheme name ----------------------
for ($r = 3; $r -le 29; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = $hklNames[$idx]
}

# --- 7. Row 2, columns C:T: the 18 "miller-index / pair" labels --------
$row2Labels = @(
    "[2, 2, 0]",
    "[2, 0, 0]",
    "[4, 0, 0]",
    "[2, 1, 1]",
    "[3, 2, 1]",
    "[2, 2, 2]",
    "[3, 1, 0]",
    "[1, 1, 0]",
    "1Pair-A",
    "1Pair-B",
    "2Pairs-A",
    "2Pairs-B",
    "3Pairs-A",
    "3Pairs-B",
    "3Pairs-C",
    "4Pairs",
    "5A4F",
    "MaxUnique"
)
for ($i = 0; $i -lt $row2Labels.Count; $i++) {
    $ws.Cells.Item(2, $i + 3).Value = $row2Labels[$i]
}

# --- 8. Rows 3..29, columns C:T = 1 -------------------------------------
for ($r = 3; $r -le 29; $r++) {
    for ($c = 3; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

# --- 9. Cosmetics --------------------------------------------------------
$ws.Range("A1").Select()
